$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header cells C1/D1 get bold+wrap style (matches existing bold header style)
$ws.Range("C1").WrapText = $true
$ws.Range("D1").WrapText = $true

# Row 3
$ws.Range("B3").Value = 'No'
$ws.Range("C3").Value = 'The requirements do state the different types of triangles to account for (5)'
$ws.Range("D3").Value = 'The requirements outline to classify triangles based on if they are equilateral, scalene, isosceles or right. Many forget to consider combined circumstances such as right isosceles or right scalene triangles.'
$ws.Range("C3").WrapText = $true
$ws.Range("D3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 51

# Row 4
$ws.Range("B4").Value = 'NA'
$ws.Range("D4").Value = 'Who is the program for?'
$ws.Range("D4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 54

# Row 5
$ws.Range("B5").Value = 'NA'
$ws.Range("D5").Value = 'What is the program being used for?'
$ws.Range("D5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 36

# Row 6
$ws.Range("B6").Value = 'Yes'
$ws.Range("C6").Value = 'Combined parameters are not mentioned in the requirements, but should have been assumed (6).'
$ws.Range("D6").Value = 'The necessities for classifying triangles are provided, and list the specifications that are being developed.'
$ws.Range("C6").WrapText = $true
$ws.Range("D6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 51

# Row 8
$ws.Range("B8").Value = 'No'
$ws.Range("C8").Value = 'The requirements provide the four classifications of triangles (5)'
$ws.Range("D8").Value = 'The requirements do not regard instances of illegal triangles, and restrictions on inputs, requiring multiple interpretations to be made by the developer.'
$ws.Range("C8").WrapText = $true
$ws.Range("D8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 36

# Row 9
$ws.Range("B9").Value = 'Yes'
$ws.Range("C9").Value = 'Right triangles can have properties of scalene or isosceles triangles (6)'
$ws.Range("D9").Value = 'The requirements should note that improper arguments should be accounted for, including negative values, or the calculation of a legal triangle.'
$ws.Range("C9").WrapText = $true
$ws.Range("D9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 36

# Row 10
$ws.Range("B10").Value = 'Yes'
$ws.Range("C10").Value = 'None'
$ws.Range("D10").Value = 'The requirements should outline an issue that must be solved.'
$ws.Range("C10").WrapText = $true
$ws.Range("D10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 54

# Row 11
$ws.Range("B11").Value = 'NA'
$ws.Range("D11").Value = 'The requirements are all that is provided.'
$ws.Range("D11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 54

# Row 13
$ws.Range("B13").Value = 'Yes'
$ws.Range("C13").Value = 'The right triangle classification can have multiple properties such as scalene or isosceles (6)'
$ws.Range("D13").Value = 'The requirements should provide the specific cases that can not be accounted for without further thought or interpretation.'
$ws.Range("C13").WrapText = $true
$ws.Range("D13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 51

# Row 14
$ws.Range("B14").Value = 'Yes'
$ws.Range("C14").Value = 'The requirements do not specify (although it should be implied) that the program must verify if the triangle is legal (3). The precision of the program should be provided (4.2), and the inputs of the program should be further specified (2).'
$ws.Range("D14").Value = 'The requirements should outline the inputs that will be tested for, accounting the verification of the triangle and the range of numbers that can be provided.'
$ws.Range("C14").WrapText = $true
$ws.Range("D14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 119

# Row 15
$ws.Range("B15").Value = 'NA'
$ws.Range("D15").Value = 'The program does perform in a timely manner, however this case would depend on the stakeholders and their use cases.'
$ws.Range("D15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 54

# Row 17
$ws.Range("B17").Value = 'Yes'
$ws.Range("C17").Value = 'The right triangle specification does not account for different types of right triangles (6).'
$ws.Range("D17").Value = 'The requirements should reference the right triangle classification more in order to guide developers to check if a triangle is a right triangle first.'
$ws.Range("C17").WrapText = $true
$ws.Range("D17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 51

# Row 18
$ws.Range("B18").Value = 'Yes'
$ws.Range("C18").Value = 'The requirements do not specified the different types of right triangles (6).'
$ws.Range("D18").Value = 'While the requirements do specify what triangle classifications to account for, the requirements do not identify the specific circumstances of triangle classification.'
$ws.Range("C18").WrapText = $true
$ws.Range("D18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 51

# Row 19
$ws.Range("B19").Value = 'NA'
$ws.Range("D19").Value = 'Use cases must be identified for the stakeholders in order to scale the consistensy of requirements.'
$ws.Range("D19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 54

# Row 21
$ws.Range("B21").Value = 'Yes'
$ws.Range("C21").Value = 'The requirements do not specified the different types of right triangles (6).'
$ws.Range("D21").Value = 'The direct intent of the system is provided, as the triangle classifications needed are provided in the requirements.'
$ws.Range("C21").WrapText = $true
$ws.Range("D21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 51

# Row 22
$ws.Range("B22").Value = 'Yes'
$ws.Range("D22").Value = 'There should be more specifications in the requirements regarding inputs into the program. Specifically, the developer needs to be provided valid inputs (2), and that the triangle should be a legal triangle (3).'
$ws.Range("D22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 51

# Row 23
$ws.Range("B23").Value = 'Yes'
$ws.Range("C23").Value = 'The requirements do not clearly state that a right triangle is separately verified as it can share properties of other triangles (6).'
$ws.Range("D23").Value = 'The requirements should develop more on the classification of the right triangle, instead of providing this requirement with the other classifications. Upon reading, it is hard for developers to identify whether the program should test for right triangles along with equilateral, scalene and isosceles properties.'
$ws.Range("C23").WrapText = $true
$ws.Range("D23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 68

# Final selection cell, matches the diff (activeCell D22)
$ws.Range("D22").Select() | Out-Null